# Applies the lattice-multiplication-exercises update:
#  - regenerates the multiplication problems/answers in the existing cells
#  - drops the old "34 x 23 / 72 x 44 / 53 x 86" row
#  - appends a new "21 x 74 / 49 x 55 / 46 x 86" row

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$br = [char]11

function Set-Cell($row, $col, $top, $nums, $line1, $line2) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $top + $br + $nums + $br + "  ----" + $br + $line1 + $br + $line2
}

# Row 3 (34 x 23 / 72 x 44 / 53 x 86) is removed entirely in the new version.
$t.Rows.Item(3).Delete()

# --- Row 1 ---
Set-Cell 1 1 "80 x 90" "  9    0" "8|    |" "0|    |"
Set-Cell 1 2 "64 x 31" "  3    1" "6|    |" "4|    |"
Set-Cell 1 3 "87 x 81" "  8    1" "8|    |" "7|    |"

# --- Row 2 ---
Set-Cell 2 1 "18 x 61" "  6    1" "1|    |" "8|    |"
Set-Cell 2 2 "60 x 33" "  3    3" "6|    |" "0|    |"
Set-Cell 2 3 "63 x 40" "  4    0" "6|    |" "3|    |"

# --- Row 3 (was row 4 before the delete: 38 x 47 / 55 x 50 / 20 x 55) ---
Set-Cell 3 1 "38 x 98" "  9    8" "3|    |" "8|    |"
Set-Cell 3 2 "70 x 89" "  8    9" "7|    |" "0|    |"
Set-Cell 3 3 "98 x 56" "  5    6" "9|    |" "8|    |"

# --- Row 4 (was row 5 before the delete: 74 x 91 / 66 x 47 / 13 x 43) ---
Set-Cell 4 1 "24 x 49" "  4    9" "2|    |" "4|    |"
Set-Cell 4 2 "73 x 28" "  2    8" "7|    |" "3|    |"
Set-Cell 4 3 "25 x 91" "  9    1" "2|    |" "5|    |"

# --- Row 5: brand-new row appended at the bottom ---
$t.Rows.Add() | Out-Null
$newRowIndex = $t.Rows.Count
Set-Cell $newRowIndex 1 "21 x 74" "  7    4" "2|    |" "1|    |"
Set-Cell $newRowIndex 2 "49 x 55" "  5    5" "4|    |" "9|    |"
Set-Cell $newRowIndex 3 "46 x 86" "  8    6" "4|    |" "6|    |"

Write-Output "lattice exercises updated"
